$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update quota values for rows 4,5,7,8,9,11,12 (30 -> 40)
$ws.Range("C4").Value = 40
$ws.Range("C5").Value = 40
$ws.Range("C7").Value = 40
$ws.Range("C8").Value = 40
$ws.Range("C9").Value = 40
$ws.Range("C11").Value = 40
$ws.Range("C12").Value = 40

# Add new course rows 13-16
$ws.Range("A13").Value = 12
$ws.Range("B13").Value = "Cross-platform Mobile Development with Flutter"
$ws.Range("C13").Value = 40

$ws.Range("A14").Value = 13
$ws.Range("B14").Value = "Advanced Programming in C/C++ (Russian only)"
$ws.Range("C14").Value = 40

$ws.Range("A15").Value = 14
$ws.Range("B15").Value = "Introduction to Mechanical Engineering"
$ws.Range("C15").Value = 40

$ws.Range("A16").Value = 15
$ws.Range("B16").Value = "Introduction to Electronic and Logic Circuits"
$ws.Range("C16").Value = 40

# Course-name cells all use the smaller (size 10) font, same as the rest of
# column B
$ws.Range("B13:B16").Font.Size = 10

# Highlight the last new row's course name cell with a light red/pink fill
# (FFF4CCCC => R=244 G=204 B=204 -> OLE color = R + G*256 + B*65536)
$ws.Range("B16").Interior.Color = 13421812

# Update selection to match the recorded end state
$ws.Range("P17").Select()
